$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append after the existing data (row 275 = "01-10-2021").
$dates = @("02-10-2021", "03-10-2021", "04-10-2021", "05-10-2021", "06-10-2021")
$startRow = 276

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i

    # Write the date as a literal text value (not an Excel date serial) by
    # routing it through a formula-result + paste-values round trip, which
    # avoids Excel's automatic date-string recognition on plain .Value
    # assignment and keeps the cell's style untouched (no numberFormat churn).
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Formula = "=""" + $dates[$i] + """"
    $cellA.Copy()
    $cellA.PasteSpecial(-4163)

    $ws.Cells.Item($row, 2).Value = 3068
    $ws.Cells.Item($row, 3).Value = 204
}

$excel.CutCopyMode = $false
